$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.289.12"
$ws.Range("E2").Value = "  +1.76%  "

$ws.Range("D3").Value = "3.164.60"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "3.160.01"
$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.69%  "

$ws.Range("E10").Value = "  -1.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.79%  "

$ws.Range("D15").Value = "3.687.68"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "3.168.90"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("D18").Value = "63.330.26"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.693"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "

$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").Value = "0.0₃0733"
$ws.Range("E37").Value = "  +5.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0389"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("E41").Value = "  -2.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "389.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.50%  "

$ws.Range("D44").Value = "2.785.25"
$ws.Range("E44").Value = "  -7.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "127.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.47%  "

